$wb = $excel.ActiveWorkbook
$ov = $wb.Worksheets.Item("Overview")
$mc = $wb.Worksheets.Item("Microcode")

# ---------------------------------------------------------------------------
# Microcode sheet: fill in the CPX/CPY microcode rows (rows 54-57) that were
# previously blank placeholders, and remove the stray B58 cell.
# ---------------------------------------------------------------------------

# Helper style "donor" cells that already carry the three cell styles we
# need to reproduce (11 = filled microcode cell, 12 = trailing empty/unused
# cycle cell, 13 = highlighted opcode-bit cell).
$style11 = $mc.Range("C53")
$style12 = $mc.Range("M53")
$style13 = $mc.Range("G25")

function Set-CellStyle($range, $donor) {
    $donor.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

# Row 54 - CPX #
$mc.Range("B54").Value = "CPX #"
$mc.Range("F54").Value = "1111 0111 000 0 0 0 0 0"
$mc.Range("G54").Value = "1001 0110 000 0 0 0 0 0"
$mc.Range("H54").Value = "0000 0000 001 1 1 0 1 1"
foreach ($col in @("I","J","K","L","M","N","O","P","Q","R")) {
    $mc.Range("$col" + "54").Value = "0000 0000 000 0 0 0 0 0"
}
Set-CellStyle $mc.Range("B54") $style11
Set-CellStyle $mc.Range("F54") $style11
Set-CellStyle $mc.Range("G54") $style13
Set-CellStyle $mc.Range("H54") $style11
Set-CellStyle $mc.Range("I54:R54") $style12

# Row 55 - CPY #
$mc.Range("B55").Value = "CPY #"
$mc.Range("F55").Value = "1111 0111 000 0 0 0 0 0"
$mc.Range("G55").Value = "1010 0110 000 0 0 0 0 0"
$mc.Range("H55").Value = "0000 0000 001 1 1 0 1 1"
foreach ($col in @("I","J","K","L","M","N","O","P","Q","R")) {
    $mc.Range("$col" + "55").Value = "0000 0000 000 0 0 0 0 0"
}
Set-CellStyle $mc.Range("B55") $style11
Set-CellStyle $mc.Range("F55") $style11
Set-CellStyle $mc.Range("G55") $style13
Set-CellStyle $mc.Range("H55") $style11
Set-CellStyle $mc.Range("I55:R55") $style12

# Row 56 - CPX abs
$mc.Range("B56").Value = "CPX abs"
$mc.Range("F56").Value = "1111 0111 000 0 0 0 0 1"
$mc.Range("G56").Value = "1111 0001 000 0 1 0 0 0"
$mc.Range("H56").Value = "0111 0001 000 0 0 0 0 0"
$mc.Range("I56").Value = "1111 0111 000 0 0 0 0 1"
$mc.Range("J56").Value = "1001 0110 000 0 0 0 0 0"
$mc.Range("K56").Value = "0000 0000 001 1 1 0 1 0"
foreach ($col in @("L","M","N","O","P","Q","R")) {
    $mc.Range("$col" + "56").Value = "0000 0000 000 0 0 0 0 0"
}
Set-CellStyle $mc.Range("B56") $style11
Set-CellStyle $mc.Range("F56") $style11
Set-CellStyle $mc.Range("G56") $style11
Set-CellStyle $mc.Range("H56") $style11
Set-CellStyle $mc.Range("I56") $style11
Set-CellStyle $mc.Range("J56") $style11
Set-CellStyle $mc.Range("K56") $style11
Set-CellStyle $mc.Range("L56:R56") $style12

# Row 57 - CPY abs
$mc.Range("B57").Value = "CPY abs"
$mc.Range("F57").Value = "1111 0111 000 0 0 0 0 1"
$mc.Range("G57").Value = "1111 0001 000 0 1 0 0 0"
$mc.Range("H57").Value = "0111 0001 000 0 0 0 0 0"
$mc.Range("I57").Value = "1111 0111 000 0 0 0 0 1"
$mc.Range("J57").Value = "1010 0110 000 0 0 0 0 0"
$mc.Range("K57").Value = "0000 0000 001 1 1 0 1 0"
foreach ($col in @("L","M","N","O","P","Q","R")) {
    $mc.Range("$col" + "57").Value = "0000 0000 000 0 0 0 0 0"
}
Set-CellStyle $mc.Range("B57") $style11
Set-CellStyle $mc.Range("F57") $style11
Set-CellStyle $mc.Range("G57") $style11
Set-CellStyle $mc.Range("H57") $style11
Set-CellStyle $mc.Range("I57") $style11
Set-CellStyle $mc.Range("J57") $style11
Set-CellStyle $mc.Range("K57") $style11
Set-CellStyle $mc.Range("L57:R57") $style12

# Row 58 - remove the now-unused, blank B58 cell entirely
$mc.Range("B58").Clear()

# ---------------------------------------------------------------------------
# Restore the cursor / selection state recorded in the saved workbook.
# ---------------------------------------------------------------------------
$ov.Range("C18").Select() | Out-Null

$mc.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$mc.Range("F58").Select() | Out-Null
